$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 25; this shifts old rows 25-49 down to 26-50
# and grows the used range from A1:AA49 to A1:AA50.
$ws.Rows("25:25").Insert()

# Populate the newly inserted row 25 with the new event record
# ("Concerto della Banda cittadina A.Ferri").
$ws.Range("A25").Value = "Concerti,Musica"
$ws.Range("B25").Value = "Modena"
$ws.Range("C25").Value = "Piazza Roma"
$ws.Range("D25").Value = "2022-05-30T09:17:47+00:00"
$ws.Range("E25").Value = "In occasione della Festa della Repubblica"
$ws.Range("F25").Value = "2022-05-30T09:20:03+00:00"
$ws.Range("G25").Value = ""
$ws.Range("H25").Value = "2022-06-02T09:00:00+00:00"
$ws.Range("I25").Value = "2022-06-02T10:00:00+00:00"
$ws.Range("J25").Value = "https://www.comune.modena.it/api/novita/eventi/2022/concerto-della-banda-cittadina-a-ferri/@@images/eaeaa0eb-f871-4bc3-885d-dd19e1c86208.jpeg"
$ws.Range("K25").Value = "Palazzo Ducale sede dell'Accademia Militare"
$ws.Range("L25").Value = "2022-05-30T09:20:03+00:00"
$ws.Range("M25").Value = "Cortile d'onore del Palazzo Ducale - Accedemia Miliatre"
$ws.Range("N25").Value = " ore 21.00"
$ws.Range("O25").Value = ""
$ws.Range("P25").Value = " ingresso libero"
$ws.Range("Q25").Value = ""
$ws.Range("R25").Value = ""
$ws.Range("S25").Value = "Concerto della Banda cittadina A.Ferri"
$ws.Range("T25").Value = ""
$ws.Range("U25").Value = ""
$ws.Range("V25").Value = $false
$ws.Range("W25").Value = 41123
$ws.Range("X25").Value = "https://www.comune.modena.it/novita/eventi/2022/concerto-della-banda-cittadina-a-ferri"
$ws.Range("Y25").Value = "44,64582"
$ws.Range("Z25").Value = "10,92572"
$ws.Range("AA25").Value = "POINT (10.92572 44.64582)"
